$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from
# 2023-10-07 (45206) to 2023-10-09 (45208) for every data row (2..172).
$ws.Range("C2:C172").Value = 45208
